# Charcrete database: convert comma-decimal "dot style" text amounts in
# column D back into real numeric values so downstream formulas (e.g. the
# LN(D..) uncertainty columns) evaluate instead of raising #VALUE!.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

$ws.Range("D13").Value = 0.15739
$ws.Range("D14").Value = 0.047136
$ws.Range("D36").Value = 0.01
$ws.Range("D37").Value = 0.108
$ws.Range("D38").Value = 0.198
$ws.Range("D49").Value = 2.7
$ws.Range("D50").Value = 2.7
$ws.Range("D51").Value = 2.7
$ws.Range("D62").Value = 18.3
$ws.Range("D64").Value = 5.4
$ws.Range("D65").Value = 0.185
$ws.Range("D67").Value = 4.37
$ws.Range("D68").Value = 0.2
$ws.Range("D70").Value = 325.6

# scroll the view so column B is the left-most visible column
$excel.ActiveWindow.ScrollColumn = 2
